# Adjust layout infos in body e-mail
# Re-populate the "DataBase" roster table with the next turn's staff/date
# data (24/10, 25/10, 26/10) and shrink the table from 10 to 7 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataBase")

# --- 1. Write the new header + data rows (A1:E7) -----------------------
$ws.Range("A1").Value = "Nome"
$ws.Range("B1").Value = "Data"
$ws.Range("C1").Value = "Codigo"
$ws.Range("D1").Value = "Unidade"
$ws.Range("E1").Value = "Matrícula"

$ws.Range("A2").Value = "NIVALDO GALVAO DE OLIVEIRA"
$ws.Range("B2").Value = "24/10"
$ws.Range("C2").Value = 4462
$ws.Range("D2").Value = "AG SÃO SEBASTIAO"
$ws.Range("E2").Value = "c150713;"

$ws.Range("A3").Value = "GUILHERME MARTINS DOS SANTOS JUNIOR"
$ws.Range("B3").Value = "24/10"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "AG LUZIANIA"
$ws.Range("E3").Value = "c150714;"

$ws.Range("A4").Value = "JOAREZ DE MOURA"
$ws.Range("B4").Value = "24/10"
$ws.Range("C4").Value = 688
$ws.Range("D4").Value = "AG GAMA"
$ws.Range("E4").Value = "c150715;"

$ws.Range("A5").Value = "ARTHUR DE CASTRO"
$ws.Range("B5").Value = "25/10"
$ws.Range("C5").Value = 4462
$ws.Range("D5").Value = "AG SÃO SEBASTIAO"
$ws.Range("E5").Value = "c150716;"

$ws.Range("A6").Value = "WILTON VASQUEZ"
$ws.Range("B6").Value = "26/10"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "AG LUZIANIA"
$ws.Range("E6").Value = "c150717;"

$ws.Range("A7").Value = "AILTON MARCELO"
$ws.Range("B7").Value = "26/10"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "AG LUZIANIA"
$ws.Range("E7").Value = "c150717;"

# --- 2. Drop the now-unused trailing rows (old rows 8:10) --------------
$ws.Range("A8:E10").Clear()

# --- 3. Shrink the table / ListObject to the new extent -----------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E7"))

# Remove the leftover sort state from the previous (larger) table range.
$lo.Sort.SortFields.Clear()

# The table was effectively rebuilt for this turn, so it gets a fresh name.
$lo.Name = "Tabela13"

# --- 4. Leave the selection where the user left off (first empty row) ---
$ws.Range("A8").Select()
